# "M01 Froze Decoder 3" - refresh the per-epoch accuracy numbers in column B
# (new training run output) and update the address text in A102:A109, then
# leave the selection where the author left off (F93, having scrolled the
# sheet down so row 81 is at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> new column-B value (only rows whose value actually changed)
$bUpdates = @{
    5   = 0.9375
    6   = 0.921875
    8   = 0.890625
    9   = 0.90625
    11  = 0.875
    14  = 0.90625
    16  = 0.921875
    17  = 0.8125
    18  = 0.828125
    19  = 0.828125
    20  = 0.84375
    21  = 0.828125
    22  = 0.8125
    23  = 0.84375
    24  = 0.796875
    25  = 0.859375
    26  = 0.828125
    27  = 0.8125
    28  = 0.8125
    29  = 0.84375
    31  = 0.828125
    32  = 0.84375
    33  = 0.84375
    34  = 0.859375
    36  = 0.859375
    37  = 0.84375
    38  = 0.828125
    39  = 0.84375
    40  = 0.8125
    41  = 0.828125
    42  = 0.796875
    43  = 0.84375
    44  = 0.796875
    46  = 0.78125
    47  = 0.8125
    48  = 0.78125
    49  = 0.796875
    50  = 0.75
    51  = 0.78125
    52  = 0.765625
    53  = 0.75
    54  = 0.78125
    55  = 0.78125
    56  = 0.78125
    57  = 0.765625
    58  = 0.765625
    59  = 0.765625
    60  = 0.765625
    61  = 0.765625
    62  = 0.765625
    63  = 0.765625
    64  = 0.765625
    65  = 0.765625
    66  = 0.765625
    67  = 0.765625
    68  = 0.765625
    69  = 0.765625
    70  = 0.765625
    71  = 0.78125
    72  = 0.78125
    73  = 0.78125
    74  = 0.78125
    75  = 0.78125
    76  = 0.78125
    77  = 0.796875
    78  = 0.796875
    79  = 0.796875
    80  = 0.796875
    81  = 0.796875
    82  = 0.796875
    83  = 0.796875
    84  = 0.796875
    85  = 0.796875
    86  = 0.796875
    87  = 0.796875
    88  = 0.765625
    89  = 0.765625
    90  = 0.78125
    91  = 0.78125
    92  = 0.78125
    93  = 0.78125
    94  = 0.78125
    95  = 0.78125
    96  = 0.78125
    97  = 0.78125
    98  = 0.78125
    99  = 0.78125
    100 = 0.796875
    101 = 0.796875
    102 = 0.796875
    103 = 0.75
    104 = 0.75
    105 = 0.6875
    108 = 0.609375
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# column A, rows 102-109: the repr of the (rerun) DisplayOutputs object
# changed its in-memory address
$newRepr = "<__main__.DisplayOutputs object at 0x7f6604d9e310>"
for ($row = 102; $row -le 109; $row++) {
    $ws.Cells.Item($row, 1).Value = $newRepr
}

# Scroll the window so row 81 is the top visible row, then leave the
# selection on F93 (where the author ended up).
$excel.ActiveWindow.ScrollRow = 81
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F93").Select()
